# "added 4wk low sales check"
# Updates the MyForecast (D), Inventory Coverage (H), Stockout Risk (I),
# Reorder Urgency (J) and Seasonality Index (L) columns on the
# "Forecast Comparison" sheet, and refreshes the dependent roll-up
# statistics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$fc = $wb.Worksheets.Item("Forecast Comparison")
$sm = $wb.Worksheets.Item("Summary")

# New MyForecast values per row (2..17)
$forecast = @{
    2  = 9
    3  = 9
    4  = 9
    5  = 9
    6  = 9
    7  = 9
    8  = 9
    9  = 9
    10 = 8
    11 = 8
    12 = 8
    13 = 8
    14 = 8
    15 = 8
    16 = 7
    17 = 7
}

# New Inventory Coverage values per row (only rows that change)
$coverage = @{
    2  = 6.67
    3  = 5.67
    4  = 4.67
    5  = 3.67
    6  = 2.67
    7  = 1.67
    8  = 0.67
    9  = 0
    10 = 0
}

# New Seasonality Index values per row (2..17)
$seasonality = @{
    2  = 0.95
    3  = 0.86
    4  = 0.84
    5  = 0.97
    6  = 1.06
    7  = 0.8
    8  = 0.97
    9  = 0.88
    10 = 1.02
    11 = 1.06
    12 = 1.03
    13 = 0.9
    14 = 0.86
    15 = 1.14
    16 = 0.97
    17 = 0.97
}

foreach ($row in $forecast.Keys) {
    $fc.Range("D$row").Value = $forecast[$row]
}

foreach ($row in $coverage.Keys) {
    $fc.Range("H$row").Value = $coverage[$row]
}

foreach ($row in $seasonality.Keys) {
    $fc.Range("L$row").Value = $seasonality[$row]
}

# Stock dropped low enough to flip risk/urgency flags on rows 8 and 9
$fc.Range("I9").Value = "High"
$fc.Range("J8").Value = "Urgent"
$fc.Range("J9").Value = "Urgent"

# Refresh dependent Summary sheet roll-ups.
# These cells hold their numbers as text (matching the rest of the
# Summary column), so force a text format before assigning the value.
$summaryUpdates = @{
    "B9"  = "134"
    "B10" = "72"
    "B11" = "36"
    "B12" = "9"
    "B14" = "7"
}
foreach ($addr in $summaryUpdates.Keys) {
    $cell = $sm.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
}
